$d = $word.ActiveDocument

# 1. Update the first paragraph's text to add trailing spaces, then append
#    a red-colored parenthetical note as three separate runs.
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Microsoft word document.  ", 2)

$para = $d.Paragraphs(1)
$r = $para.Range
$tail = $d.Range($r.Start, $r.End - 1)
$tail.Collapse(0)
$tail.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$tail.Font.Color = 255
$tail.Collapse(0)
$tail.InsertAfter("rsion for main branch")
$tail.Font.Color = 255
$tail.Collapse(0)
$tail.InsertAfter(")")
$tail.Font.Color = 255

# 2. Remove the trailing "ank God almighty, we are free at last." paragraph
#    entirely (including its paragraph mark).
$count = $d.Paragraphs.Count
$last = $d.Paragraphs($count)
$last.Range.Delete()
